# Update report metadata header
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header fields (rows 5, 8, 9) ---
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:47 PM"
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 18

# --- Insert the new rows needed (bottom-most block first so earlier, still-
#     untouched row numbers above stay valid) ---
# Wednesday block gains 3 extra line-item rows (6 -> 9) before its TOTAL row (37)
$ws.Range("A37:A39").EntireRow.Insert()
# Monday block gains 1 extra line-item row (2 -> 3) before its TOTAL row (18)
$ws.Range("A18").EntireRow.Insert()

# Helper to write one data-row's worth of cells (columns A-H; G is left blank)
function Set-LineRow {
    param($row, $point, $code, $worktype, $desc, $uom, $units, $price)
    $ws.Cells.Item($row, 1).Value = $point
    $ws.Cells.Item($row, 2).Value = $code
    $ws.Cells.Item($row, 3).Value = $worktype
    $ws.Cells.Item($row, 4).Value = $desc
    $ws.Cells.Item($row, 5).Value = $uom
    $ws.Cells.Item($row, 6).Value = $units
    $ws.Cells.Item($row, 8).Value = $price
}

function Set-TotalRow {
    param($row, $total)
    $ws.Cells.Item($row, 1).Value = "TOTAL"
    $ws.Cells.Item($row, 8).Value = $total
}

# --- Monday (06/30/2025) : rows 16-18 data, row 19 TOTAL ---
Set-LineRow 16 "Point 03" "CNA-TR"   "Inst" "CNA,Transfer Conductor"                  "EA" 3  0
Set-LineRow 17 "Point 03" "POL-40-4" "Inst" "Pole,40ft,Class 4"                       "EA" 1  0
Set-LineRow 18 "Point 03" "PLA-BACK" "Inst" "Difficult Location Equip Adder-Backyard" "EA" 18 0
Set-TotalRow 19 0

# Row 18 was a brand-new row inserted with the formatting of row 17 (the
# "even" stripe); restripe it to match row 16's ("odd" stripe) formatting,
# matching the banded look of the other data rows.
$ws.Range("A16:H16").Copy()
$ws.Range("A18:H18").PasteSpecial(-4122)
$excel.CutCopyMode = $false
Set-LineRow 18 "Point 03" "PLA-BACK" "Inst" "Difficult Location Equip Adder-Backyard" "EA" 18 0

# --- Tuesday (07/01/2025) : rows 24-26 data, row 27 TOTAL (unchanged counts,
#     only pricing collapses to 0) ---
Set-LineRow 24 "Point 05" "POL-40-4" "Inst" "Pole,40ft,Class 4"                       "EA" 1 0
Set-LineRow 25 "Point 05" "PLA-BACK" "Inst" "Difficult Location Equip Adder-Backyard" "EA" 8 0
Set-LineRow 26 "Point 05" "PLA-BACK" "Inst" "Difficult Location Equip Adder-Backyard" "EA" 8 0
Set-TotalRow 27 0

# --- Wednesday (07/02/2025) : rows 32-40 data, row 41 TOTAL ---
Set-LineRow 32 "Point 07" "GYD-MPY"         "Inst" "GYD,Marker-Plastic-Yellow"              "EA" 1  0
Set-LineRow 33 "Point 07" "GYF-38-42W-I"    "Inst" "GYF,3/8,42In Wire Mt,Insulator Adder"   "EA" 1  0
Set-LineRow 34 "Point 07" "GYF-38-D-78P-EP" "Inst" "GYF,3/8,Down,78in Pole mt,EyePlate"     "EA" 1  0
Set-LineRow 35 "Point 07" "GYW-38"          "Inst" "GYW,3/8 in. EHS (15,400 lbs)"           "FT" 40 0
Set-LineRow 36 "Point 07" "POL-40-4"        "Inst" "Pole,40ft,Class 4"                       "EA" 1  0
Set-LineRow 37 "Point 07" "PLA-BACK"        "Inst" "Difficult Location Equip Adder-Backyard" "EA" 18 0
Set-LineRow 38 "Point 07" "PLA-BACK"        "Inst" "Difficult Location Equip Adder-Backyard" "EA" 18 0
Set-LineRow 39 "Point 07" "PLA-HDIG"        "Inst" "PLA,Hand Dig or Additional  Excavation"  "EA" 1  0
Set-LineRow 40 "Point 08" "ANC-SHM-10-84-D1" "Rem" "ANC,Sgl Hlx Mach,10in,84in,Db Eye 1in"   "EA" 1  0
Set-TotalRow 41 0

# Rows 38 and 40 were brand-new inserted rows that picked up the "even"
# stripe formatting; restripe them to the "odd" stripe like rows 32/34/36.
$ws.Range("A32:H32").Copy()
$ws.Range("A38:H38").PasteSpecial(-4122)
$ws.Range("A40:H40").PasteSpecial(-4122)
$excel.CutCopyMode = $false
Set-LineRow 38 "Point 07" "PLA-BACK"         "Inst" "Difficult Location Equip Adder-Backyard" "EA" 18 0
Set-LineRow 40 "Point 08" "ANC-SHM-10-84-D1" "Rem"  "ANC,Sgl Hlx Mach,10in,84in,Db Eye 1in"    "EA" 1  0

# --- Thursday (07/03/2025) : rows 46-48 data, row 49 TOTAL (unchanged
#     counts, only pricing collapses to 0) ---
Set-LineRow 46 "Point 09" "POL-40-2" "Inst" "Pole,40ft,Class 2"                       "EA" 1  0
Set-LineRow 47 "Point 09" "PLA-BACK" "Inst" "Difficult Location Equip Adder-Backyard" "EA" 18 0
Set-LineRow 48 "Point 09" "PLA-BACK" "Inst" "Difficult Location Equip Adder-Backyard" "EA" 18 0
Set-TotalRow 49 0

# Re-merge an existing column-I-spanning region so the sheet's recorded
# dimension keeps extending through column I (to I49) after all the row
# insertions above, matching the workbook's original extent convention.
$ws.Range("G13:I13").Merge()
